# Add team record (Wins/Losses/Ties) columns to the DET_1996 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF with the same header style (s="1")
# as the existing header cells. Copy format from an existing header cell first,
# then overwrite the value/text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# Data rows 2-55: constant team record values for every player row.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 53
    $ws.Cells.Item($r, 31).Value = 109
    $ws.Cells.Item($r, 32).Value = 0
}
